{"js": "// Correct the dataset size mentioned in the introduction:\n//   \" The dataset used will contain information on 5,110 individuals (n = 5,110) \"\n// becomes\n//   \" The dataset used contains information on 4,909 individuals (n = 4,909) \"\nconst body = context.document.body;\n\nconst target = body.search(\n  \"The dataset used will contain information on 5,110 individuals (n = 5,110)\",\n  { matchCase: true }\n);\ntarget.load(\"items/text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\n    \"The dataset used contains information on 4,909 individuals (n = 4,909)\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n} else {\n  // Fallback: perform the two edits independently if the full phrase\n  // could not be located as a single match (e.g. already partially edited).\n  const verbPhrase = body.search(\"will contain\", { matchCase: true });\n  verbPhrase.load(\"items\");\n  await context.sync();\n  if (verbPhrase.items.length > 0) {\n    verbPhrase.items[0].insertText(\"contains\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  let numbers = body.search(\"5,110\", { matchCase: true });\n  numbers.load(\"items\");\n  await context.sync();\n  while (numbers.items.length > 0) {\n    numbers.items[0].insertText(\"4,909\", Word.InsertLocation.replace);\n    await context.sync();\n    numbers = body.search(\"5,110\", { matchCase: true });\n    numbers.load(\"items\");\n    await context.sync();\n  }\n}\n", "ps1": "# Correct the dataset size mentioned in the introduction:\n#   \" The dataset used will contain information on 5,110 individuals (n = 5,110) \"\n# becomes\n#   \" The dataset used contains information on 4,909 individuals (n = 4,909) \"\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"The dataset used will contain information on 5,110 individuals (n = 5,110)\")\n\nif ($found) {\n  $rng.Text = \"The dataset used contains information on 4,909 individuals (n = 4,909)\"\n} else {\n  # Fallback: perform the two edits independently if the full phrase\n  # could not be located as a single match.\n  $verb = $d.Content\n  if ($verb.Find.Execute(\"will contain\")) {\n    $verb.Text = \"contains\"\n  }\n\n  $num1 = $d.Content\n  while ($num1.Find.Execute(\"5,110\")) {\n    $num1.Text = \"4,909\"\n    $num1 = $d.Content\n  }\n}\n"}
